$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks (this is the only reliable way the runtime
# supports removing hyperlinks - per-item deletion is a no-op here).
$ws.Hyperlinks.Delete()

# New set of usernames/emails (omkarhundre+215 .. omkarhundre+230) that
# replace the old two rows and extend the sheet to 16 data rows.
$startNum = 215
$endNum = 230

$firstDataRow = 2
$row = $firstDataRow
$testNo = 1
for ($n = $startNum; $n -le $endNum; $n++) {
    $email = "omkarhundre+$n@arcitech.ai"
    $ws.Range("A$row").Value2 = $testNo
    $ws.Range("B$row").Value2 = $email
    $ws.Hyperlinks.Add($ws.Range("B$row"), "mailto:$email") | Out-Null
    $row = $row + 1
    $testNo = $testNo + 1
}
$lastHyperlinkRow = $row - 1

# Two trailing rows that only carry the Testcase_no. counter (no username).
$ws.Range("A$row").Value2 = $testNo
$row = $row + 1
$testNo = $testNo + 1
$ws.Range("A$row").Value2 = $testNo

$lastDataRow = $row

# Re-apply the built-in "Hyperlink" style to the whole used B column so the
# mail cells keep looking like the original two (underlined hyperlink font)
# instead of the ad-hoc style variant created internally by Hyperlinks.Add.
$ws.Range("B$($firstDataRow):B$lastHyperlinkRow").Style = "Hyperlink"

# Match the original workbook's behaviour of leaving the selection on the
# first empty cell below the data.
$ws.Range("B$lastDataRow").Select() | Out-Null
